# Daily auto-push update: a new reading for 2026/01/21 (19:00, rank 201) was
# appended to the source log. Because the sheet is sorted chronologically and
# 2026/01/21 sorts before the already-present 2026/12/29 block, the new
# record lands as row 687, pushing every following row down by one
# (687:728 -> 688:729) and extending the used range to D729.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 687:728 down to 688:729, opening up a blank row 687.
$ws.Rows("687:687").Insert()

# Give the new A687 cell the same "plain text" formatting the rest of the
# date column uses, so the date string isn't auto-converted to a date serial.
$ws.Range("A687").NumberFormat = "@"
$ws.Range("A687").Value = "2026/01/21"

# Re-apply the neighboring (unstyled) date cell's formatting so A687 ends up
# with the same default style as every other data-row cell in column A.
$ws.Range("A686").Copy()
$ws.Range("A687").PasteSpecial(-4122)

$ws.Range("B687").Value = "水"
$ws.Range("C687").Value = 19
$ws.Range("D687").Value = 201
